$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap dates between row 2 and row 4 (column D)
$d2 = $ws.Range("D2").Value2
$d4 = $ws.Range("D4").Value2
$ws.Range("D2").Value2 = $d4
$ws.Range("D4").Value2 = $d2

# Swap volumes between row 2 and row 4 (column M)
$m2 = $ws.Range("M2").Value2
$m4 = $ws.Range("M4").Value2
$ws.Range("M2").Value2 = $m4
$ws.Range("M4").Value2 = $m2
